# Apply the documented edits to the resource-alternate.xlsx example workbook.
$wb = $excel.ActiveWorkbook

# --- "resource" worksheet: cost of "Carl" (row 4) changes from 100 to 120 ---
$wsResource = $wb.Worksheets.Item("resource")
$wsResource.Activate()
$wsResource.Range("E4").Value = 120
$wsResource.Range("F7").Select() | Out-Null

# --- "parameter" worksheet: add a new "plan.webservice" parameter row ---
$wsParameter = $wb.Worksheets.Item("parameter")
$wsParameter.Activate()
$wsParameter.Range("A15").Value = "plan.webservice"

# Writing the literal text "true" directly would be auto-coerced to a
# Boolean by the COM layer (same as typing TRUE into a cell in Excel),
# but the existing sheet stores this value as the shared text string
# "true" (it is reused elsewhere on this sheet). Round-trip through a
# formula + paste-values so the result lands back in the cell as
# literal text instead of a boolean.
$wsParameter.Range("B15").Formula = "=""true"""
$wsParameter.Range("B15").Copy() | Out-Null
$wsParameter.Range("B15").PasteSpecial(-4163) | Out-Null
